$wb = $excel.ActiveWorkbook

# --- "Hidden Row & Col" sheet: add two new "notes" plus real hidden data ---
$ws = $wb.Worksheets.Item("Hidden Row & Col")

# Make room: 4 new rows (3..6) and 2 new columns (C..D)
$ws.Rows("3:6").Insert()
$ws.Columns("C:D").Insert()

# Row 2 / Col B+C now hold real (non-blank) numeric data instead of the old note
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("B3").Value = 2

# The three explanatory notes, now placed in column D (always visible)
$ws.Range("D4").Value = "Column A or Row 1 are hidden with no data"
$ws.Range("D5").Value = "Column B and Row 2 are hidden with data"
$ws.Range("D6").Value = "Column C and Row 3 have zero width/height"

# Hide row 2 and row 3 (row 1 was already hidden), and columns A, B, C
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Columns.Item(1).Hidden = $true
$ws.Columns.Item(2).Hidden = $true
$ws.Columns.Item(3).Hidden = $true

# This sheet becomes the active one, selection moves to A4
$ws.Activate()
$ws.Range("A4").Select()

Write-Host "done"
